# Automatische test-sync: 2025-08-05 19:31:50
#
# Adds a new test-mail log row (row 46) to the "Logs" sheet, extends the
# conditional-formatting ranges to include it, and re-syncs the "Dashboard"
# pivot-style summary (the "Overig" category now outranks
# "Retour / Terugbetaling" since it grew from 3 to 4).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "Logs" sheet: append the new row with the 5th test mail.
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$newRow = 46

$logs.Cells.Item($newRow, 1).Value  = "Kun jij deze klant even bellen vandaag?"
$logs.Cells.Item($newRow, 2).Value  = "mailmind.test@zohomail.eu"
$logs.Cells.Item($newRow, 3).Value  = "Testmail #5: Kun jij deze klant even bellen vandaag?"
$logs.Cells.Item($newRow, 4).Value  = "Overig"
$logs.Cells.Item($newRow, 5).Value  = "Bedankt, we hebben dit doorgestuurd naar support@bedrijf.nl."
$logs.Cells.Item($newRow, 6).Value  = "2025-08-05 19:30:54"
$logs.Cells.Item($newRow, 7).Value  = "Ja"
$logs.Cells.Item($newRow, 8).Value  = "Ja"
$logs.Cells.Item($newRow, 9).Value  = "Nee"
$logs.Cells.Item($newRow, 10).Value = "Nee"

# ---------------------------------------------------------------------
# 2. Grow the conditional-formatting ranges (D/G/H/I/J, rows 2..45) so
#    they also cover the freshly-added row 46.
# ---------------------------------------------------------------------
$columns = @("D", "G", "H", "I", "J")
foreach ($col in $columns) {
    $oldRange = $logs.Range($col + "2:" + $col + "45")
    $newRange = $logs.Range($col + "2:" + $col + "46")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# ---------------------------------------------------------------------
# 3. "Dashboard" sheet: re-sync the category summary. "Overig" grew from
#    3 to 4 mentions and now outranks "Retour / Terugbetaling" (still 3),
#    so the two rows swap places in the sorted-by-count table.
# ---------------------------------------------------------------------
$dashboard = $wb.Worksheets.Item("Dashboard")

$dashboard.Cells.Item(5, 1).Value = "Overig"
$dashboard.Cells.Item(5, 2).Value = 4

$dashboard.Cells.Item(6, 1).Value = "Retour / Terugbetaling"
$dashboard.Cells.Item(6, 2).Value = 3
